$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting D:K to E:L
$ws.Columns("D:D").Insert()

# Copy formatting from the (now-shifted) old column D -- now column E -- into new column D
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

# Populate the new column D with the updated period values
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 125400
$ws.Range("D9").Value = 121700
$ws.Range("D10").Value = 3700
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 13300
$ws.Range("D15").Value = 1900
$ws.Range("D17").Value = 150400
$ws.Range("D18").Value = -25100
$ws.Range("D20").Value = 2400
$ws.Range("D21").Value = -13500
$ws.Range("D22").Value = 1500
$ws.Range("D23").Value = -24200
$ws.Range("D24").Value = -200
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -24000
$ws.Range("D27").Value = -24000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = -100
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -2400
$ws.Range("D33").Value = -24100
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -24100
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 1200
$ws.Range("D42").Value = "NA"
$ws.Range("D43").Value = 17500
$ws.Range("D44").Value = 22700
$ws.Range("D45").Value = 1800
$ws.Range("D46").Value = 43100
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 49100
$ws.Range("D49").Value = 6600
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 400
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 99200
$ws.Range("D57").Value = 11600
$ws.Range("D58").Value = 12900
$ws.Range("D59").Value = 27300
$ws.Range("D60").Value = 51900
$ws.Range("D61").Value = 2000
$ws.Range("D62").Value = 2000
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 55800
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -336300
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 43400
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = -24100
$ws.Range("D83").Value = 9200
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 1900
$ws.Range("D91").Value = -2300
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -1600
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 800
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 1100
